# Update the "Comp. image meas." value for the last row (row 6), which
# drives the computed "% image meas" formula in P6 and the chart series
# that plots column P.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I6").Value = 14

# Recalculate and refresh chart caches so the chart's cached series
# values (e.g. the "% image meas" series pulling from column P) pick up
# the new computed results.
$wb.RefreshAll()
$excel.CalculateFullRebuild()

# Update the active selection on the sheet to match the author's final
# selection state (B6:K6, active cell B6).
$ws.Range("B6:K6").Select()

$wb.Save()
